# Delta method workbook restructure:
#   - Rename existing "main" sheet (the big raw-data sheet) to "main_data"
#   - Insert a new small "main" summary sheet in front of it, containing one
#     summarised row (mean coefficient / mean years) for the Aslam (2007) study
#   - Copy over the three header comments (B1/C1/F1) onto the new sheet
#   - Update the "main_data" sheetView (no longer tab-selected, frozen pane /
#     selection now parked on row 166)

$wb = $excel.ActiveWorkbook

# 1. Rename the current "main" sheet (580-row raw dataset) to "main_data"
$mainData = $wb.Worksheets.Item("main")
$mainData.Name = "main_data"

# 2. Insert a brand-new sheet in front of it and call it "main"
$main = $wb.Worksheets.Add()
$main.Name = "main"
$main.Move($mainData)

# 3. Populate the new summary sheet
$main.Range("A1").Value = "study"
$main.Range("B1").Value = "coef_mean"
$main.Range("C1").Value = "years_mean"
$main.Range("D1").Value = "beta_higher"
$main.Range("E1").Value = "se_higher"
$main.Range("F1").Value = "beta_lower"
$main.Range("G1").Value = "se_lower"

$main.Range("A2").Value = "Aslam (2007)"
$main.Range("B2").Value = 5.9
$main.Range("C2").Value = 12
$main.Range("D2").Value = 1.57
$main.Range("E2").Value = 0.878
$main.Range("F2").Value = -0.867
$main.Range("G2").Value = 0.197

# Styling carried over from the template: bold/shaded header (s=1), bold/
# shaded + centered header for the two new columns (s=2), centered percent
# column (s=3), centered integer column (s=6)
$main.Range("A1").Style = $mainData.Range("A1").Style
$main.Range("D1").Style = $mainData.Range("A1").Style
$main.Range("E1").Style = $mainData.Range("A1").Style
$main.Range("F1").Style = $mainData.Range("A1").Style
$main.Range("G1").Style = $mainData.Range("A1").Style
$main.Range("B1").Style = $mainData.Range("B1").Style
$main.Range("C1").Style = $mainData.Range("B1").Style
$main.Range("B2").Style = $mainData.Range("B2").Style
$main.Range("C2").Style = $mainData.Range("C166").Style

$main.Range("C2").Selection

# 4. Re-create the header comments on the new sheet
$main.Range("B1").AddComment("Petr Čala:`nMake sure to get this right")
$main.Range("C1").AddComment("Petr Čala:`nalways the number of years as a difference from the reference category, not as a whole")
$main.Range("F1").AddComment("Petr Čala:`nLeave these empty for the single coefficient function to be used")

# 5. Update the main_data sheet view: no longer the selected tab, pane/
#    selection now resting on row 166 (Aslam (2007), the last of the rows
#    with a C166 comment-anchored 14-year category)
$mainData.Activate()
$mainData.Range("A166:G166").Select()
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.FreezePanes = $true

$wb.Worksheets.Item("main").Select()
